$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Build a list of target cell references covered by this edit so we can force
# them to text formatting (the source data stores every column, including the
# numeric-looking ones, as text) before writing the new values.
$targetCells = @(
    "B138", "D138", "E138", "H138", "B139", "D139", "E139", "H139", "B140", "D140", "E140", "H140", "B141", "D141", "E141", "H141", "B142", "D142", "E142", "H142", "B143", "D143", "E143", "H143", "B144", "D144", "E144", "H144", "B145", "D145", "E145", "H145", "D146", "E146", "H146", "B147", "D147", "E147", "H147", "B148", "D148", "E148", "H148", "B149", "D149", "E149", "H149", "B150", "D150", "E150", "H150", "B151", "D151", "E151", "H151", "B152", "D152", "E152", "H152", "B153", "D153", "E153", "H153", "B154", "D154", "E154", "F154", "G154", "H154", "B155", "D155", "E155", "F155", "G155", "H155", "B156", "D156", "E156", "F156", "G156", "H156", "B157", "D157", "E157", "F157", "G157", "H157", "B158", "D158", "E158", "F158", "G158", "H158", "B159", "D159", "E159", "F159", "G159", "H159", "D160", "E160", "F160", "G160", "H160", "B161", "D161", "E161", "F161", "G161", "H161", "B162", "D162", "E162", "F162", "G162", "H162", "B163", "D163", "E163", "F163", "G163", "H163", "B164", "D164", "E164", "F164", "G164", "H164", "B165", "D165", "E165", "F165", "G165", "H165", "B166", "D166", "E166", "F166", "G166", "H166", "B167", "D167", "E167", "F167", "G167", "H167", "B168", "D168", "E168", "F168", "G168", "H168", "B169", "D169", "E169", "F169", "G169", "H169"
)

foreach ($ref in $targetCells) {
    $ws.Range($ref).NumberFormat = "@"
}

$ws.Range("B138").Value = "s137_e104_24019-32_1_8.jpeg"
$ws.Range("D138").Value = "780"
$ws.Range("E138").Value = "978"
$ws.Range("H138").Value = "18"

$ws.Range("B139").Value = "s138_e102_24019-32_1_6.jpeg"
$ws.Range("D139").Value = "1008"
$ws.Range("E139").Value = "1192"
$ws.Range("H139").Value = "5"

$ws.Range("B140").Value = "s139_e116_24019-32_3_4.jpeg"
$ws.Range("D140").Value = "2488"
$ws.Range("E140").Value = "869"
$ws.Range("H140").Value = "110"

$ws.Range("B141").Value = "s140_e120_24019-32_3_8.jpeg"
$ws.Range("D141").Value = "1496"
$ws.Range("E141").Value = "1798"
$ws.Range("H141").Value = "71"

$ws.Range("B142").Value = "s141_e109_24019-32_2_3.jpeg"
$ws.Range("D142").Value = "2106"
$ws.Range("E142").Value = "1111"
$ws.Range("H142").Value = "145"

$ws.Range("B143").Value = "s142_e115_24019-32_3_3.jpeg"
$ws.Range("D143").Value = "1315"
$ws.Range("E143").Value = "1684"
$ws.Range("H143").Value = "115"

$ws.Range("B144").Value = "s143_e123_24019-32_4_5.jpeg"
$ws.Range("D144").Value = "1369"
$ws.Range("E144").Value = "1235"
$ws.Range("H144").Value = "153"

$ws.Range("B145").Value = "s144_e110_24019-32_2_2.jpeg"
$ws.Range("D145").Value = "1724"
$ws.Range("E145").Value = "660"
$ws.Range("H145").Value = "153"

$ws.Range("D146").Value = "1654"
$ws.Range("E146").Value = "376"
$ws.Range("H146").Value = "48"

$ws.Range("B147").Value = "s146_e117_24019-32_3_5.jpeg"
$ws.Range("D147").Value = "1433"
$ws.Range("E147").Value = "1747"
$ws.Range("H147").Value = "176"

$ws.Range("B148").Value = "s147_e106_24019-32_2_6.jpeg"
$ws.Range("D148").Value = "98"
$ws.Range("E148").Value = "1548"
$ws.Range("H148").Value = "105"

$ws.Range("B149").Value = "s148_e119_24019-32_3_7.jpeg"
$ws.Range("D149").Value = "2550"
$ws.Range("E149").Value = "1704"
$ws.Range("H149").Value = "119"

$ws.Range("B150").Value = "s149_e100_24019-32_1_4.jpeg"
$ws.Range("D150").Value = "2253"
$ws.Range("E150").Value = "430"
$ws.Range("H150").Value = "116"

$ws.Range("B151").Value = "s150_e121_24019-32_4_7.jpeg"
$ws.Range("D151").Value = "1273"
$ws.Range("E151").Value = "1730"
$ws.Range("H151").Value = "49"

$ws.Range("B152").Value = "s151_e114_24019-32_3_2.jpeg"
$ws.Range("D152").Value = "1223"
$ws.Range("E152").Value = "1020"
$ws.Range("H152").Value = "147"

$ws.Range("B153").Value = "s152_e105_24019-32_2_7.jpeg"
$ws.Range("D153").Value = "2387"
$ws.Range("E153").Value = "1929"
$ws.Range("H153").Value = "86"

$ws.Range("B154").Value = "s153_e141_24019-32_3_1.jpeg"
$ws.Range("D154").Value = "1669"
$ws.Range("E154").Value = "714"
$ws.Range("F154").Value = "20"
$ws.Range("G154").Value = "20"
$ws.Range("H154").Value = "46"

$ws.Range("B155").Value = "s154_e127_24019-32_1_3.jpeg"
$ws.Range("D155").Value = "168"
$ws.Range("E155").Value = "164"
$ws.Range("F155").Value = "40"
$ws.Range("G155").Value = "20"
$ws.Range("H155").Value = "30"

$ws.Range("B156").Value = "s155_e139_24019-32_2_1.jpeg"
$ws.Range("D156").Value = "1160"
$ws.Range("E156").Value = "1253"
$ws.Range("F156").Value = "61"
$ws.Range("G156").Value = "20"
$ws.Range("H156").Value = "45"

$ws.Range("B157").Value = "s156_e143_24019-32_3_3.jpeg"
$ws.Range("D157").Value = "1544"
$ws.Range("E157").Value = "224"
$ws.Range("F157").Value = "81"
$ws.Range("G157").Value = "20"
$ws.Range("H157").Value = "20"

$ws.Range("B158").Value = "s157_e126_24019-32_1_2.jpeg"
$ws.Range("D158").Value = "1730"
$ws.Range("E158").Value = "1319"
$ws.Range("F158").Value = "102"
$ws.Range("G158").Value = "20"
$ws.Range("H158").Value = "110"

$ws.Range("B159").Value = "s158_e150_24019-32_4_6.jpeg"
$ws.Range("D159").Value = "450"
$ws.Range("E159").Value = "1140"
$ws.Range("F159").Value = "122"
$ws.Range("G159").Value = "20"
$ws.Range("H159").Value = "167"

$ws.Range("D160").Value = "1842"
$ws.Range("E160").Value = "716"
$ws.Range("F160").Value = "143"
$ws.Range("G160").Value = "20"
$ws.Range("H160").Value = "95"

$ws.Range("B161").Value = "s160_e132_24019-32_1_8.jpeg"
$ws.Range("D161").Value = "1454"
$ws.Range("E161").Value = "1514"
$ws.Range("F161").Value = "163"
$ws.Range("G161").Value = "20"
$ws.Range("H161").Value = "93"

$ws.Range("B162").Value = "s161_e147_24019-32_3_7.jpeg"
$ws.Range("D162").Value = "1079"
$ws.Range("E162").Value = "654"
$ws.Range("F162").Value = "184"
$ws.Range("G162").Value = "20"
$ws.Range("H162").Value = "156"

$ws.Range("B163").Value = "s162_e136_24019-32_2_4.jpeg"
$ws.Range("D163").Value = "734"
$ws.Range("E163").Value = "305"
$ws.Range("F163").Value = "204"
$ws.Range("G163").Value = "20"
$ws.Range("H163").Value = "33"

$ws.Range("B164").Value = "s163_e137_24019-32_2_3.jpeg"
$ws.Range("D164").Value = "873"
$ws.Range("E164").Value = "1149"
$ws.Range("F164").Value = "30"
$ws.Range("G164").Value = "30"
$ws.Range("H164").Value = "126"

$ws.Range("B165").Value = "s164_e140_24019-32_2_0.jpeg"
$ws.Range("D165").Value = "116"
$ws.Range("E165").Value = "516"
$ws.Range("F165").Value = "61"
$ws.Range("G165").Value = "30"
$ws.Range("H165").Value = "89"

$ws.Range("B166").Value = "s165_e144_24019-32_3_4.jpeg"
$ws.Range("D166").Value = "1487"
$ws.Range("E166").Value = "88"
$ws.Range("F166").Value = "92"
$ws.Range("G166").Value = "30"
$ws.Range("H166").Value = "118"

$ws.Range("B167").Value = "s166_e138_24019-32_2_2.jpeg"
$ws.Range("D167").Value = "868"
$ws.Range("E167").Value = "381"
$ws.Range("F167").Value = "122"
$ws.Range("G167").Value = "30"
$ws.Range("H167").Value = "60"

$ws.Range("B168").Value = "s167_e146_24019-32_3_6.jpeg"
$ws.Range("D168").Value = "1547"
$ws.Range("E168").Value = "1483"
$ws.Range("F168").Value = "153"
$ws.Range("G168").Value = "30"
$ws.Range("H168").Value = "25"

$ws.Range("B169").Value = "s168_e148_24019-32_3_8.jpeg"
$ws.Range("D169").Value = "606"
$ws.Range("E169").Value = "1070"
$ws.Range("F169").Value = "184"
$ws.Range("G169").Value = "30"
$ws.Range("H169").Value = "26"

foreach ($ref in $targetCells) {
    $ws.Range($ref).Style = "Normal"
}
